$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 5 values and add the missing "PASS" in F5
$ws.Range("C5").Value = "xeumeaauflczvbx@gmail.com"
$ws.Range("D5").Value = "mluklKTDCO5"
$ws.Range("E5").Value = "pass"
$ws.Range("F5").Value = "PASS"

# Add new row 6
$ws.Range("C6").Value = "vhnjnquirwqrosv@gmail.com"
$ws.Range("D6").Value = "yicovWOHEI5"
$ws.Range("E6").Value = "pass"
$ws.Range("F6").Value = "PASS"

# Add new row 7
$ws.Range("C7").Value = "ubizvtrlheyavla@gmail.com"
$ws.Range("D7").Value = "ipnyyYETFZ5"
$ws.Range("E7").Value = "pass"
$ws.Range("F7").Value = "PASS"

# Add new row 8
$ws.Range("C8").Value = "lddawxqfztqllxu@gmail.com"
$ws.Range("D8").Value = "yqdesHMCPL5"
$ws.Range("E8").Value = "pass"
$ws.Range("F8").Value = "PASS"

# Update the selected cell to match the new active selection
$ws.Range("D13").Select()
